# Daily update at 8 AM UTC
# Appends the next day's row (row 55) to the "Wins Over Time" tracker and
# moves the "latest day" date-only formatting from the old last row (54)
# down to the new last row (55).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 54 ("A54") was formatted with the special "latest row" date-only
# format (YYYY-MM-DD). Now that it's no longer the last row, give it the
# same format as the rest of the date column (YYYY-MM-DD HH:MM:SS), matching
# row 53's format.
$ws.Range("A54").NumberFormat = $ws.Range("A53").NumberFormat

# Append the new day's data as row 55.
$ws.Cells.Item(55, 1).Value = 45640
$ws.Cells.Item(55, 2).Value = 133
$ws.Cells.Item(55, 3).Value = 120
$ws.Cells.Item(55, 4).Value = 123

# The new last row takes on the "latest row" date-only format.
$ws.Range("A55").NumberFormat = "YYYY-MM-DD"
